$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 28-44 use the same centered style (s=2) as the rest of the data table;
# set alignment first so freshly-created cells (rows 31+) pick it up too.
$ws.Range("A28:C44").HorizontalAlignment = -4108

# Fill in newly added proverb rows (28-41): batal / meaning / source
$ws.Range("A28").Value = "بارگہ ءِ گرند ءُ چلگ ءِ گروگ۔"
$ws.Range("B28").Value = "گرند ءُ جمبر کہ بیت گڑا ہئور الم بیت، ہئور کہ بیت بارگہ انت"
$ws.Range("C28").Value = "لال ءُ یاقوت"

$ws.Range("A29").Value = "بازارءَ سوداے کن انت  پٹاٹہ و پیماز، بہاے نہ کن انت غیرت و میار۔"
$ws.Range("B29").Value = "غیرت ءُ لج ہنچیں چیزاَنت کہ پہ بہانہ اَنت"
$ws.Range("C29").Value = "لال ءُ یاقوت"

$ws.Range("A30").Value = "بان پہ بانک، مال پہ مالک۔"
$ws.Range("B30").Value = "کمالی مرد آباد بیت مالدار بیت، بودیں جنین گس ءَ آباد کنت"
$ws.Range("C30").Value = "لال ءُ یاقوت"

$ws.Range("A31").Value = "بخت ءَ و عزت ءَ کپگ ءَ میل۔"
$ws.Range("B31").Value = "روزی کہ کیت بسم الله کن، عزات ءَ بہ دار"
$ws.Range("C31").Value = "لال ءُ یاقوت"

$ws.Range("A32").Value = "بخت پہ دوچن زورت نہ بیت۔"
$ws.Range("B32").Value = "بخت چہ انسان ءِ جوڑ بوہگ ءَ گوں آئی گون اِنت"
$ws.Range("C32").Value = "لال ءُ یاقوت"

$ws.Range("A33").Value = "بخت ءُ اقبال، نہ پہ اصل انت نہ پہ کمال۔"
$ws.Range("B33").Value = "بخت پہ خدائی نیمگ ءَ، مال اؤں خدائی  نیمگ ءَ، کس دلرنج مہ بیت، کمال ءُ جمال نہ اِنت"
$ws.Range("C33").Value = "لال ءُ یاقوت"

$ws.Range("A34").Value = "بخت ہر مردم ءَ گون نہ بیت۔"
$ws.Range("B34").Value = "عقل ہر کجاترا نصیب ءَ نہ بیت"
$ws.Range("C34").Value = "لال ءُ یاقوت"

$ws.Range("A35").Value = "بخت، پہ قسمت۔"
$ws.Range("B35").Value = "کئی کہ بخت کجام وڑ اِنت آئی ءَ ہما رسیت"
$ws.Range("C35").Value = "لال ءُ یاقوت"

$ws.Range("A36").Value = "بختاور نہ بیت زہر آور۔"
$ws.Range("B36").Value = "نیک ءُ پارسائیں مردم ءَ راپہ کسے ءَ حسد نئیت"
$ws.Range("C36").Value = "لال ءُ یاقوت"

$ws.Range("A37").Value = "بد عمل ءِ قسمت ٹَلی۔"
$ws.Range("B37").Value = "کسے کہ خرابیں کار کنت، خراب سو چیت گڑا وشی نہ گندیت"
$ws.Range("C37").Value = "لال ءُ یاقوت"

$ws.Range("B38").Value = "بدیں مردم ءَ را نہ بیت ہمراہ"
$ws.Range("A38").Value = "بدکار ءِ کسے نہ بیت دوستدار۔"
$ws.Range("C38").Value = "لال ءُ یاقوت"

$ws.Range("A39").Value = "بدعمل ءِ کاربے عمل بیت۔"
$ws.Range("B39").Value = "بد عملیں انسان ءِ ہچیز ءِ تہا برکت نہ بیت"
$ws.Range("C39").Value = "لال ءُ یاقوت"

$ws.Range("B40").Value = "بدکار ءُ بدیں مردم کسی یارنہ بنت"
$ws.Range("A40").Value = "بدکاری کسی یار نہ بیت۔"
$ws.Range("C40").Value = "لال ءُ یاقوت"

$ws.Range("A41").Value = "بدی، خیرنہ گندی۔"
$ws.Range("B41").Value = "آکہ بدیں مردمے آخیر نہ گندی"
$ws.Range("C41").Value = "لال ءُ یاقوت"

# View state as last saved in the source workbook: scrolled so row 17 is
# the top visible row, with B41 as the active selection.
$ws.Range("B41").Select()
try {
    $excel.ActiveWindow.ScrollRow = 17
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
    # scrolling API not available in this host; selection above still applies
}